$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) "Corgi working on recording videos for interview #2" -> "... #3"
#    The trailing "2" must become "3" while the text stays split the
#    same way it already was (" #" run followed by a distinct "3" run).
# ------------------------------------------------------------------
$p1 = $d.Paragraphs.Item(1)
$pStart = $p1.Range.Start
$pText = $p1.Range.Text
$twoOffset = $pText.IndexOf("#2") + 1
$twoStart = $pStart + $twoOffset
$twoRange = $d.Range($twoStart, $twoStart + 1)

# Toggle a formatting property on/off around the text edit so the "2"
# remains its own run instead of being re-merged with its " #" sibling.
$twoRange.Bold = 1
$textRange = $d.Range($twoStart, $twoStart + 1)
$textRange.Text = "3"
$finalRange = $d.Range($twoStart, $twoStart + 1)
$finalRange.Bold = 0

# ------------------------------------------------------------------
# 2) Add a new, empty sub-bullet paragraph (ilvl 1) right after the
#    "Helps companies to get the reputation they deserve" bullet.
# ------------------------------------------------------------------
$deserveParagraph = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $candidate = $d.Paragraphs.Item($i)
    if ($candidate.Range.Text -like "*Helps companies to get the reputation they deserve*") {
        $deserveParagraph = $candidate
        break
    }
}

$insertionPoint = $deserveParagraph.Range.End
$insertionRange = $d.Range($insertionPoint, $insertionPoint)
$insertionRange.InsertAfter([char]13)
